# Auto-generated Excel COM-interop script
# Adds MonthlyEvents and OneTimeEvents worksheets with content pass / event data

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Create MonthlyEvents sheet (sheet2) ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "MonthlyEvents"

# --- Create OneTimeEvents sheet (sheet3) ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "OneTimeEvents"

# === Populate MonthlyEvents (sheet2) ===
$ws2.Range("B1:B5").NumberFormat = "mm-dd-yy"

$ws2.Range("A1").Value = "name"
$ws2.Range("B1").Value = "start_time"
$ws2.Range("C1").Value = "duration"
$ws2.Range("D1").Value = "reccurence"
$ws2.Range("E1").Value = "icon_url"

$ws2.Range("A2").Value = "lost_city_of_gold"
$ws2.Range("B2").Value = 43791
$ws2.Range("C2").Value = "4 days"
$ws2.Range("D2").Value = "28 days"
$ws2.Range("E2").Value = "events/lost_city_of_gold"

$ws2.Range("A3").Value = "dragon_invasion"
$ws2.Range("B3").Value = 44742
$ws2.Range("C3").Value = "5 days"
$ws2.Range("D3").Value = "28 days"
$ws2.Range("E3").Value = "events/dragon_invasion"

$ws2.Range("A4").Value = "kings_caprice"
$ws2.Range("B4").Value = 44083
$ws2.Range("C4").Value = "7 days"
$ws2.Range("D4").Value = "28 days"
$ws2.Range("E4").Value = "events/kings_caprice"

$ws2.Range("A5").Value = "tower_of_titans"
$ws2.Range("B5").Value = 44008
$ws2.Range("C5").Value = "4 days"
$ws2.Range("D5").Value = "28 days"
$ws2.Range("E5").Value = "events/tower_of_titans"

# Column widths (MonthlyEvents)
$ws2.Columns.Item(1).ColumnWidth = 19.571428571428573   # -> stored ~20.33203125 (target 20.33203125)
$ws2.Columns.Item(2).ColumnWidth = 18.285714285714285   # -> stored 19 (target 19)
$ws2.Columns.Item(5).ColumnWidth = 23.142857142857142   # -> stored ~23.83203125 (target 23.83203125)

$ws2.Range("F5").Select()

# === Populate OneTimeEvents (sheet3) ===
$ws3.Range("D1:E12").NumberFormat = "mm-dd-yy"

$ws3.Range("A1").Value = "name"
$ws3.Range("B1").Value = "name_en"
$ws3.Range("C1").Value = "name_zh"
$ws3.Range("D1").Value = "start_time"
$ws3.Range("E1").Value = "end_time"
$ws3.Range("F1").Value = "icon_url"

$ws3.Range("A2").Value = "full_bloom_pass"
$ws3.Range("B2").Value = "Full Bloom Content Pass"
$ws3.Range("C2").Value = "春暖花开内容令状"
$ws3.Range("D2").Value = 45370
$ws3.Range("E2").Value = 45385
$ws3.Range("F2").Value = "events/content_pass/full_bloom"

$ws3.Range("A3").Value = "year_of_dragon_pass"
$ws3.Range("B3").Value = "Year of the Dragon Content Pass"
$ws3.Range("C3").Value = "龙年"
$ws3.Range("D3").Value = 45314
$ws3.Range("E3").Value = 45341
$ws3.Range("F3").Value = "events/content_pass/year_of_dragon"

$ws3.Range("A4").Value = "christmas_pass"
$ws3.Range("B4").Value = "Christmas Content Pass"
$ws3.Range("C4").Value = "圣诞内容令状"
$ws3.Range("D4").Value = 45282
$ws3.Range("E4").Value = 45299
$ws3.Range("F4").Value = "events/content_pass/christmas"

$ws3.Range("A5").Value = "vampire_masquerade_pass"
$ws3.Range("B5").Value = "Vampire Masquerade Content Pass"
$ws3.Range("C5").Value = "小淘气内容令状"
$ws3.Range("D5").Value = 45202
$ws3.Range("E5").Value = 45230
$ws3.Range("F5").Value = "events/content_pass/vampire"

$ws3.Range("A6").Value = "bjorn_pass"
$ws3.Range("B6").Value = "Bjron Content Pass"
$ws3.Range("C6").Value = "比约恩内容令状"
$ws3.Range("D6").Value = 45146
$ws3.Range("E6").Value = 45161
$ws3.Range("F6").Value = "events/content_pass/bjorn"

$ws3.Range("A7").Value = "spacefarer_pass"
$ws3.Range("B7").Value = "Spacefarer Content Pass"
$ws3.Range("C7").Value = "遨游太空内容令状"
$ws3.Range("D7").Value = 45062
$ws3.Range("E7").Value = 45089
$ws3.Range("F7").Value = "events/content_pass/spacefarer"

$ws3.Range("A8").Value = "year_of_the_rabbit_pass"
$ws3.Range("B8").Value = "Year of the Rabbit Content Pass"
$ws3.Range("C8").Value = "兔年"
$ws3.Range("D8").Value = 44936
$ws3.Range("E8").Value = 44965
$ws3.Range("F8").Value = "events/content_pass/year_of_rabbit"

$ws3.Range("A9").Value = "halloween2022_pass"
$ws3.Range("B9").Value = "Halloween Content Pass"
$ws3.Range("C9").Value = "2022年万圣节"
$ws3.Range("D9").Value = 44859
$ws3.Range("E9").Value = 44866
$ws3.Range("F9").Value = "events/content_pass/halloween2022"

$ws3.Range("A10").Value = "year_of_the_tiger_pass"
$ws3.Range("B10").Value = "Year of the Tiger Content Pass"
$ws3.Range("C10").Value = "虎年"
$ws3.Range("D10").Value = 44585
$ws3.Range("E10").Value = 44598
$ws3.Range("F10").Value = "events/content_pass/year_of_tiger"

$ws3.Range("A11").Value = "avatar_pass"
$ws3.Range("B11").Value = "Avatar Content Pass"
$ws3.Range("C11").Value = "降世神通：最后的气宗"
$ws3.Range("D11").Value = 44501
$ws3.Range("E11").Value = 44530
$ws3.Range("F11").Value = "events/content_pass/avatar"

$ws3.Range("A12").Value = "cinco_de_mayo_pass"
$ws3.Range("B12").Value = "Cinco de Mayo Content Pass"
$ws3.Range("C12").Value = "2022年五月五日节"
$ws3.Range("D12").Value = 44683
$ws3.Range("E12").Value = 44710
$ws3.Range("F12").Value = "events/content_pass/cinco_de_mayo"

# Column widths (OneTimeEvents)
$ws3.Columns.Item(1).ColumnWidth = 22.428571428571427   # -> stored ~23.1640625 (target 23.1640625)
$ws3.Columns.Item(2).ColumnWidth = 31.428571428571427   # -> stored ~32.1640625 (target 32.1640625)
$ws3.Columns.Item(3).ColumnWidth = 22.428571428571427   # -> stored ~23.1640625 (target 23.1640625)
$ws3.Columns.Item(4).ColumnWidth = 20.428571428571427   # -> stored ~21.1640625 (target 21.1640625)
$ws3.Columns.Item(5).ColumnWidth = 24.142857142857142   # -> stored ~24.83203125 (target 24.83203125)

$ws3.Range("F13").Select()
$ws3.Activate()

